$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 90, pushing the previous
# rows 90-97 down to 91-98 (dimension grows from R97 to R98).
$ws.Rows.Item(90).EntireRow.Insert()

# Fill the new row 90 with the new "Locoto" price record.
$ws.Range("A90").Value = 10
$ws.Range("B90").Value = "Vega Modelo de Temuco"
$ws.Range("C90").Value = "La Araucanía"
$ws.Range("D90").Value = 45166
$ws.Range("E90").Value = 9
$ws.Range("F90").Value = 100112042
$ws.Range("G90").Value = "Locoto"
$ws.Range("H90").Value = "Sin especificar"
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 250
$ws.Range("K90").Value = 1600
$ws.Range("L90").Value = 1600
$ws.Range("M90").Value = 1600
$ws.Range("N90").Value = "$/kilo"
$ws.Range("O90").Value = "Región de Arica y Parinacota"
$ws.Range("P90").Value = 1600
$ws.Range("Q90").Value = 1
$ws.Range("R90").Value = "Hortaliza"
